$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knight+Player")

# --- Sprites section: rename assets, clear old assignment/status marks ---
$ws.Range("B4").Value = "Body"
$ws.Range("B5").Value = "Legs"
$ws.Range("B6").Value = "Helmet"
$ws.Range("C4:G7").ClearContents()

# --- Script section ---
$ws.Range("B14").Value = "Script scene change "
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()

# --- Sound section: rename existing two rows ---
$ws.Range("B18").Value = "Fly damage sound"
$ws.Range("B19").Value = "Boss fly damage sound"

# Insert 8 new formatted rows (copying B19's format) right after row 19 so the
# "Prefabs" row below gets pushed from row 21 down to row 28, then close the
# one-row gap this creates immediately above "Prefabs".
$ws.Range("B19").Copy()
$ws.Range("B20:B27").Insert(-4121)
$ws.Range("A28:G28").Delete(-4162)

$ws.Range("B20").Value = "Ambient sounds of start/end screen"
$ws.Range("B21").Value = "In game music"
$ws.Range("B22").Value = "Armor clink/Take damage/Menu button"
$ws.Range("B23").Value = "Jump whoosh"
$ws.Range("B24").Value = "Hitting dummy wooden click"
$ws.Range("B25").Value = "Windmill creaking"
$ws.Range("B26").Value = "Door open"
$ws.Range("B27").Value = "Death oof"

# Match the saved selection state
$ws.Range("B18").Select()
